$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (Price in column D, Volume(1h) in column E)
# Values in column D are textual (e.g. "41.851.63"), so a leading apostrophe
# forces Excel to keep them as text instead of mis-parsing them as numbers,
# and resetting Style back to Normal avoids leaving a stray text-format style.

$ws.Range("D2").Value = "'41.851.63"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.97%  "

$ws.Range("D3").Value = "'2.280.13"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.89%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").Value = "'310.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.36%  "

$ws.Range("D6").Value = "'104.87"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.80%  "

$ws.Range("D7").Value = "'0.625"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.80%  "

$ws.Range("E8").Value = "  +0.11%  "

$ws.Range("D9").Value = "'0.603"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.03%  "

$ws.Range("D10").Value = "'40.22"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.62%  "

$ws.Range("D11").Value = "'0.0907"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.43%  "

$ws.Range("D12").Value = "'8.23"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.24%  "

$ws.Range("E13").Value = "  -0.01%  "

$ws.Range("D14").Value = "'0.964"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.33%  "

$ws.Range("D15").Value = "'15.38"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.37%  "

$ws.Range("D16").Value = "'2.625.94"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.89%  "

$ws.Range("D17").Value = "'2.280.38"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.06%  "

$ws.Range("D18").Value = "'41.934.15"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.65%  "

$ws.Range("D19").Value = "'7.52"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.67%  "

$ws.Range("D20").Value = "'0.0000104"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.93%  "

$ws.Range("D21").Value = "'74.03"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.94%  "

$ws.Range("D22").Value = "'3.42"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -7.71%  "

$ws.Range("D23").Value = "'255.67"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.29%  "

$ws.Range("D24").Value = "'2.29"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.48%  "

$ws.Range("D25").Value = "'9.28"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -6.99%  "

$ws.Range("D26").Value = "'1.01"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.44%  "

$ws.Range("D27").Value = "'10.90"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.65%  "

$ws.Range("E28").Value = "  +3.21%  "

$ws.Range("D29").Value = "'22.64"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.18%  "

$ws.Range("D30").Value = "'166.25"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.55%  "

$ws.Range("D31").Value = "'35.52"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.78%  "

$ws.Range("D32").Value = "'0.0887"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.20%  "

$ws.Range("D33").Value = "'2.91"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.73%  "

$ws.Range("D34").Value = "'5.77"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.86%  "

$ws.Range("E35").Value = "  -2.04%  "

$ws.Range("D36").Value = "'0.117"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +8.10%  "

$ws.Range("D37").Value = "'4.52"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.56%  "

$ws.Range("E38").Value = "  -1.65%  "

$ws.Range("D39").Value = "'2.73"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.71%  "

$ws.Range("E40").Value = "  -4.37%  "

$ws.Range("D41").Value = "'71.84"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.81%  "

$ws.Range("D42").Value = "'97.07"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.30%  "

$ws.Range("E43").Value = "  -2.79%  "

$ws.Range("D44").Value = "'0.226"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.55%  "

$ws.Range("E45").Value = "  +0.18%  "

$ws.Range("D46").Value = "'12.19"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.74%  "

$ws.Range("D47").Value = "'111.33"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -7.98%  "

$ws.Range("D48").Value = "'8.98"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.00%  "

$ws.Range("D49").Value = "'5.29"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.17%  "

$ws.Range("D50").Value = "'73.82"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.37%  "

$ws.Range("D51").Value = "'1.561.32"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.60%  "
